$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 22 and 23 swap their Coin/Link/Price/Volume data (Avalanche <-> Uniswap) ---
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.95%  "

# --- Price (column D) and Volume(1h) (column E) updates for the remaining rows ---
# Column D holds values as text (matching the workbook's original inline-string
# cells); when the new price happens to look like a plain number ("0.519",
# "28.79", ...) the assignment is wrapped with a temporary text NumberFormat so
# Excel does not silently convert it to a floating point number, then the
# cell style is restored to Normal so no formatting changes are introduced.

$ws.Range("D2").Value = "29.690.79"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "1.607.99"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "1.838.46"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "1.610.92"
$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.564"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("D16").Value = "29.694.21"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +13.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").Value = "0.0₃0705"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.01%  "

$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").Value = "1.430.48"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("E35").Value = "  +6.53%  "

$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0171"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.553"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("E41").Value = "  +5.75%  "

$ws.Range("E42").Value = "  +4.07%  "

$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "68.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.40%  "

$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +20.51%  "

$ws.Range("E48").Value = "  +3.19%  "

$ws.Range("D49").Value = "1.747.10"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("E51").Value = "  -0.87%  "
